# Updates cryptos list figures (Price / Volume(1h) columns, and two swapped coin rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apostrophe used as a literal text-prefix so Excel stores numeric-looking values
# (e.g. "1.001") as text instead of converting them to numbers, matching the
# original inlineStr/text cells in the workbook.
$q = "'"

$ws.Range('D2').Value = '22.144.35'
$ws.Range('E2').Value = '  -1.51%  '

$ws.Range('D3').Value = '1.560.75'
$ws.Range('E3').Value = '  -1.01%  '

$ws.Range('D4').Value = $q + '1.001'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('E5').Value = '  +0.06%  '

$ws.Range('D6').Value = $q + '290.25'
$ws.Range('E6').Value = '  +0.51%  '

$ws.Range('D7').Value = $q + '0.3784'
$ws.Range('E7').Value = '  +2.68%  '

$ws.Range('D8').Value = $q + '0.3284'

$ws.Range('D9').Value = $q + '43.67'
$ws.Range('E9').Value = '  -9.22%  '

$ws.Range('D10').Value = $q + '1.137'
$ws.Range('E10').Value = '  -0.95%  '

$ws.Range('D11').Value = $q + '0.07363'
$ws.Range('E11').Value = '  -2.49%  '

$ws.Range('D12').Value = $q + '1.002'
$ws.Range('E12').Value = '  +0.05%  '

$ws.Range('D13').Value = $q + '19.93'
$ws.Range('E13').Value = '  -4.24%  '

$ws.Range('D14').Value = $q + '5.824'
$ws.Range('E14').Value = '  -2.54%  '

$ws.Range('D15').Value = $q + '6.870'
$ws.Range('E15').Value = '  -1.20%  '

$ws.Range('D16').Value = '1.560.66'
$ws.Range('E16').Value = '  -0.99%  '

$ws.Range('D17').Value = $q + '0.00001093'
$ws.Range('E17').Value = '  -2.71%  '

$ws.Range('E18').Value = '  -1.38%  '

$ws.Range('D19').Value = $q + '85.48'
$ws.Range('E19').Value = '  -2.81%  '

$ws.Range('D20').Value = $q + '6.456'
$ws.Range('E20').Value = '  +0.79%  '

$ws.Range('E21').Value = '  +0.02%  '

$ws.Range('D22').Value = $q + '16.12'
$ws.Range('E22').Value = '  -3.04%  '

$ws.Range('D23').Value = $q + '11.73'
$ws.Range('E23').Value = '  -2.37%  '

$ws.Range('D24').Value = '22.158.86'
$ws.Range('E24').Value = '  -1.39%  '

$ws.Range('D25').Value = $q + '2.272'
$ws.Range('E25').Value = '  -5.15%  '

$ws.Range('D26').Value = $q + '2.531'
$ws.Range('E26').Value = '  -4.34%  '

$ws.Range('D27').Value = $q + '150.87'
$ws.Range('E27').Value = '  -0.36%  '

$ws.Range('D28').Value = $q + '19.10'
$ws.Range('E28').Value = '  -3.01%  '

$ws.Range('D29').Value = $q + '4.863'
$ws.Range('E29').Value = '  -2.81%  '

$ws.Range('D30').Value = '1.734.84'
$ws.Range('E30').Value = '  -1.07%  '

$ws.Range('D31').Value = $q + '121.39'
$ws.Range('E31').Value = '  -3.23%  '

$ws.Range('D32').Value = $q + '1.122'
$ws.Range('E32').Value = '  +1.59%  '

$ws.Range('D33').Value = $q + '6.033'
$ws.Range('E33').Value = '  -2.01%  '

$ws.Range('D34').Value = $q + '1.876'
$ws.Range('E34').Value = '  -6.29%  '

$ws.Range('D35').Value = $q + '9.338'
$ws.Range('E35').Value = '  -5.61%  '

$ws.Range('D36').Value = $q + '0.08183'
$ws.Range('E36').Value = '  -2.35%  '

$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = $q + '5.278'
$ws.Range('E37').Value = '  -1.95%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = $q + '0.02308'
$ws.Range('E38').Value = '  -6.75%  '

$ws.Range('E39').Value = '  -3.21%  '

$ws.Range('D40').Value = $q + '0.2141'
$ws.Range('E40').Value = '  -5.06%  '

$ws.Range('D41').Value = $q + '1.238'
$ws.Range('E41').Value = '  -4.37%  '

$ws.Range('D42').Value = $q + '11.08'
$ws.Range('E42').Value = '  -3.73%  '

$ws.Range('D43').Value = $q + '1.000'
$ws.Range('E43').Value = '  +0.06%  '

$ws.Range('D44').Value = $q + '0.5987'
$ws.Range('E44').Value = '  -5.16%  '

$ws.Range('E45').Value = '  -2.76%  '

$ws.Range('D46').Value = $q + '3.759'
$ws.Range('E46').Value = '  -0.86%  '

$ws.Range('D47').Value = $q + '0.5791'
$ws.Range('E47').Value = '  -5.89%  '

$ws.Range('D48').Value = $q + '1.987'
$ws.Range('E48').Value = '  -4.15%  '

$ws.Range('D49').Value = $q + '121.07'
$ws.Range('E49').Value = '  -3.89%  '

$ws.Range('D50').Value = $q + '1.167'
$ws.Range('E50').Value = '  -4.00%  '

$ws.Range('D51').Value = $q + '0.07000'
$ws.Range('E51').Value = '  -3.36%  '
